$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.074674634519317
$ws.Cells.Item(2, 4).Value = 1.074782120527066
$ws.Cells.Item(2, 5).Value = 1.077964081760828
$ws.Cells.Item(2, 6).Value = 1.087104781879177
$ws.Cells.Item(2, 9).Value = 1.049372733764684
$ws.Cells.Item(2, 10).Value = 1.079582927464044
$ws.Cells.Item(2, 11).Value = 1.077470423401331
$ws.Cells.Item(2, 12).Value = 1.080643999204078
$ws.Cells.Item(2, 13).Value = 1.089760905514138
$ws.Cells.Item(2, 14).Value = 1.081116059050936
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.076274162652691
$ws.Cells.Item(3, 4).Value = 1.076034648658763
$ws.Cells.Item(3, 5).Value = 1.0793533705711
$ws.Cells.Item(3, 6).Value = 1.088489096058501
$ws.Cells.Item(3, 9).Value = 1.049758446728569
$ws.Cells.Item(3, 10).Value = 1.080838584652954
$ws.Cells.Item(3, 11).Value = 1.078539074064011
$ws.Cells.Item(3, 12).Value = 1.08184968754237
$ws.Cells.Item(3, 13).Value = 1.090963366646205
$ws.Cells.Item(3, 14).Value = 1.08237349941708
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.077307832173501
$ws.Cells.Item(4, 4).Value = 1.076843713677338
$ws.Cells.Item(4, 5).Value = 1.080251265120474
$ws.Cells.Item(4, 6).Value = 1.089383729580332
$ws.Cells.Item(4, 9).Value = 1.050006038464519
$ws.Cells.Item(4, 10).Value = 1.081649334146708
$ws.Cells.Item(4, 11).Value = 1.079228595322537
$ws.Cells.Item(4, 12).Value = 1.082628243605622
$ws.Cells.Item(4, 13).Value = 1.091739788229008
$ws.Cells.Item(4, 14).Value = 1.083185400268109
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.077742076190967
$ws.Cells.Item(5, 4).Value = 1.077183513954868
$ws.Cells.Item(5, 5).Value = 1.080628490502864
$ws.Cells.Item(5, 6).Value = 1.089759573767355
$ws.Cells.Item(5, 9).Value = 1.050109651707531
$ws.Cells.Item(5, 10).Value = 1.081989761780578
$ws.Cells.Item(5, 11).Value = 1.079518004335407
$ws.Cells.Item(5, 12).Value = 1.082955170134485
$ws.Cells.Item(5, 13).Value = 1.092065806195535
$ws.Cells.Item(5, 14).Value = 1.083526311348267
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.077814969695117
$ws.Cells.Item(6, 4).Value = 1.077240548688015
$ws.Cells.Item(6, 5).Value = 1.080691813855288
$ws.Cells.Item(6, 6).Value = 1.0898226646042
$ws.Cells.Item(6, 9).Value = 1.050127021075212
$ws.Cells.Item(6, 10).Value = 1.082046897130659
$ws.Cells.Item(6, 11).Value = 1.079566570223007
$ws.Cells.Item(6, 12).Value = 1.083010040512644
$ws.Cells.Item(6, 13).Value = 1.092120523359337
$ws.Cells.Item(6, 14).Value = 1.083583527837098
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.077313635773723
$ws.Cells.Item(7, 4).Value = 1.076848255397651
$ws.Cells.Item(7, 5).Value = 1.080256306597082
$ws.Cells.Item(7, 6).Value = 1.089388752641874
$ws.Cells.Item(7, 9).Value = 1.050007424810833
$ws.Cells.Item(7, 10).Value = 1.081653884566455
$ws.Cells.Item(7, 11).Value = 1.079232464243465
$ws.Cells.Item(7, 12).Value = 1.082632613493774
$ws.Cells.Item(7, 13).Value = 1.091744146021207
$ws.Cells.Item(7, 14).Value = 1.083189957149974
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.075215483357623
$ws.Cells.Item(8, 4).Value = 1.075205712433056
$ws.Cells.Item(8, 5).Value = 1.078433823368199
$ws.Cells.Item(8, 6).Value = 1.087572850898171
$ws.Cells.Item(8, 9).Value = 1.049503500690571
$ws.Cells.Item(8, 10).Value = 1.08000764769232
$ws.Cells.Item(8, 11).Value = 1.077831989366085
$ws.Cells.Item(8, 12).Value = 1.081051803019285
$ws.Cells.Item(8, 13).Value = 1.090167628167636
$ws.Cells.Item(8, 14).Value = 1.081541382430651
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.071507696774884
$ws.Cells.Item(9, 4).Value = 1.072300335015093
$ws.Cells.Item(9, 5).Value = 1.075213913047059
$ws.Cells.Item(9, 6).Value = 1.084364223443376
$ws.Cells.Item(9, 9).Value = 1.04860018061858
$ws.Cells.Item(9, 10).Value = 1.077093126048527
$ws.Cells.Item(9, 11).Value = 1.075348872514474
$ws.Cells.Item(9, 12).Value = 1.07825366138223
$ws.Cells.Item(9, 13).Value = 1.087376703105192
$ws.Cells.Item(9, 14).Value = 1.078622721831824
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.069028208677662
$ws.Cells.Item(10, 4).Value = 1.07035566366083
$ws.Cells.Item(10, 5).Value = 1.073061221939245
$ws.Cells.Item(10, 6).Value = 1.082218853728818
$ws.Cells.Item(10, 9).Value = 1.047987514620517
$ws.Cells.Item(10, 10).Value = 1.07514054620216
$ws.Cells.Item(10, 11).Value = 1.073682842990497
$ws.Cells.Item(10, 12).Value = 1.076379429180137
$ws.Cells.Item(10, 13).Value = 1.085507057550718
$ws.Cells.Item(10, 14).Value = 1.076667369098112
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.067952636669193
$ws.Cells.Item(11, 4).Value = 1.06951168137976
$ws.Cells.Item(11, 5).Value = 1.072127551311295
$ws.Cells.Item(11, 6).Value = 1.081288308179319
$ws.Cells.Item(11, 9).Value = 1.047719713164534
$ws.Cells.Item(11, 10).Value = 1.074292700640059
$ws.Cells.Item(11, 11).Value = 1.072958841081285
$ws.Cells.Item(11, 12).Value = 1.075565695708596
$ws.Cells.Item(11, 13).Value = 1.084695258945246
$ws.Cells.Item(11, 14).Value = 1.075818319498063
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.067552820171334
$ws.Cells.Item(12, 4).Value = 1.069197892397405
$ws.Cells.Item(12, 5).Value = 1.071780505216817
$ws.Cells.Item(12, 6).Value = 1.080942416298775
$ws.Cells.Item(12, 9).Value = 1.047619859482577
$ws.Cells.Item(12, 10).Value = 1.073977409964382
$ws.Cells.Item(12, 11).Value = 1.072689517413035
$ws.Cells.Item(12, 12).Value = 1.075263104418958
$ws.Cells.Item(12, 13).Value = 1.084393378830932
$ws.Cells.Item(12, 14).Value = 1.075502581073463
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.067638595958715
$ws.Cells.Item(13, 4).Value = 1.06926521472292
$ws.Cells.Item(13, 5).Value = 1.071854958733196
$ws.Cells.Item(13, 6).Value = 1.081016622526835
$ws.Cells.Item(13, 9).Value = 1.047641295689841
$ws.Cells.Item(13, 10).Value = 1.074045057438809
$ws.Cells.Item(13, 11).Value = 1.072747306329403
$ws.Cells.Item(13, 12).Value = 1.075328026527458
$ws.Cells.Item(13, 13).Value = 1.084458148727484
$ws.Cells.Item(13, 14).Value = 1.075570324615062
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.067919593915142
$ws.Cells.Item(14, 4).Value = 1.069485749570936
$ws.Cells.Item(14, 5).Value = 1.072098869298324
$ws.Cells.Item(14, 6).Value = 1.081259721707081
$ws.Cells.Item(14, 9).Value = 1.047711467001927
$ws.Cells.Item(14, 10).Value = 1.074266646065875
$ws.Cells.Item(14, 11).Value = 1.072936586841725
$ws.Cells.Item(14, 12).Value = 1.075540690282331
$ws.Cells.Item(14, 13).Value = 1.08467031245738
$ws.Cells.Item(14, 14).Value = 1.075792227923395
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.068092685819414
$ws.Cells.Item(15, 4).Value = 1.069621588887762
$ws.Cells.Item(15, 5).Value = 1.072249118712798
$ws.Cells.Item(15, 6).Value = 1.08140947033476
$ws.Cells.Item(15, 9).Value = 1.047754651387743
$ws.Cells.Item(15, 10).Value = 1.074403125789541
$ws.Cells.Item(15, 11).Value = 1.073053155994946
$ws.Cells.Item(15, 12).Value = 1.075671674908416
$ws.Cells.Item(15, 13).Value = 1.08480098801017
$ws.Cells.Item(15, 14).Value = 1.075928901463923
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.069099549171921
$ws.Cells.Item(16, 4).Value = 1.070411634802077
$ws.Cells.Item(16, 5).Value = 1.073123153456587
$ws.Cells.Item(16, 6).Value = 1.082280576911267
$ws.Cells.Item(16, 9).Value = 1.04800523454085
$ws.Cells.Item(16, 10).Value = 1.075196764425221
$ws.Cells.Item(16, 11).Value = 1.073730837278844
$ws.Cells.Item(16, 12).Value = 1.076433387482517
$ws.Cells.Item(16, 13).Value = 1.085560886378509
$ws.Cells.Item(16, 14).Value = 1.076723667157498
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.069730602169731
$ws.Cells.Item(17, 4).Value = 1.070906689022527
$ws.Cells.Item(17, 5).Value = 1.073670994320917
$ws.Cells.Item(17, 6).Value = 1.082826569072149
$ws.Cells.Item(17, 9).Value = 1.048161743992433
$ws.Cells.Item(17, 10).Value = 1.075693954231966
$ws.Cells.Item(17, 11).Value = 1.07415522815443
$ws.Cells.Item(17, 12).Value = 1.076910601368141
$ws.Cells.Item(17, 13).Value = 1.08603694867727
$ws.Cells.Item(17, 14).Value = 1.077221563030796
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.070098498472374
$ws.Cells.Item(18, 4).Value = 1.071195260340483
$ws.Cells.Item(18, 5).Value = 1.073990392451892
$ws.Cells.Item(18, 6).Value = 1.083144884517394
$ws.Cells.Item(18, 9).Value = 1.048252791022804
$ws.Cells.Item(18, 10).Value = 1.075983728910264
$ws.Cells.Item(18, 11).Value = 1.074402517738184
$ws.Cells.Item(18, 12).Value = 1.07718874244026
$ws.Cells.Item(18, 13).Value = 1.086314413098443
$ws.Cells.Item(18, 14).Value = 1.077511749222375
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.070223910320085
$ws.Cells.Item(19, 4).Value = 1.071293624457854
$ws.Cells.Item(19, 5).Value = 1.074099274059263
$ws.Cells.Item(19, 6).Value = 1.083253396209876
$ws.Cells.Item(19, 9).Value = 1.048283794694065
$ws.Cells.Item(19, 10).Value = 1.076082496192286
$ws.Cells.Item(19, 11).Value = 1.07448679481287
$ws.Cells.Item(19, 12).Value = 1.077283546001983
$ws.Cells.Item(19, 13).Value = 1.086408985094969
$ws.Cells.Item(19, 14).Value = 1.077610656765267
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.069662915491619
$ws.Cells.Item(20, 4).Value = 1.070853593620165
$ws.Cells.Item(20, 5).Value = 1.073612231563359
$ws.Cells.Item(20, 6).Value = 1.082768005113354
$ws.Cells.Item(20, 9).Value = 1.048144977090014
$ws.Cells.Item(20, 10).Value = 1.075640634081298
$ws.Cells.Item(20, 11).Value = 1.074109720980422
$ws.Cells.Item(20, 12).Value = 1.07685942257796
$ws.Cells.Item(20, 13).Value = 1.085985893952219
$ws.Cells.Item(20, 14).Value = 1.077168167159398
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.067836855415806
$ws.Cells.Item(21, 4).Value = 1.069420815792345
$ws.Cells.Item(21, 5).Value = 1.072027050382921
$ws.Cells.Item(21, 6).Value = 1.081188141904303
$ws.Cells.Item(21, 9).Value = 1.047690813819693
$ws.Cells.Item(21, 10).Value = 1.074201403823
$ws.Cells.Item(21, 11).Value = 1.072880859463396
$ws.Cells.Item(21, 12).Value = 1.075478075397036
$ws.Cells.Item(21, 13).Value = 1.084607845023206
$ws.Cells.Item(21, 14).Value = 1.075726893029052
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.066686991467204
$ws.Cells.Item(22, 4).Value = 1.068518253249452
$ws.Cells.Item(22, 5).Value = 1.07102899596174
$ws.Cells.Item(22, 6).Value = 1.080193393136993
$ws.Cells.Item(22, 9).Value = 1.047403061688218
$ws.Cells.Item(22, 10).Value = 1.073294398100653
$ws.Cells.Item(22, 11).Value = 1.072105924574606
$ws.Cells.Item(22, 12).Value = 1.074607629014241
$ws.Cells.Item(22, 13).Value = 1.083739428957716
$ws.Cells.Item(22, 14).Value = 1.074818599254544
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.067296725138491
$ws.Cells.Item(23, 4).Value = 1.068996883922856
$ws.Cells.Item(23, 5).Value = 1.07155821765809
$ws.Cells.Item(23, 6).Value = 1.080720865974971
$ws.Cells.Item(23, 9).Value = 1.047555814150318
$ws.Cells.Item(23, 10).Value = 1.073775421059178
$ws.Cells.Item(23, 11).Value = 1.072516952495354
$ws.Cells.Item(23, 12).Value = 1.07506925526512
$ws.Cells.Item(23, 13).Value = 1.084199982958393
$ws.Cells.Item(23, 14).Value = 1.075300305320845
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.069693500763578
$ws.Cells.Item(24, 4).Value = 1.070877585720258
$ws.Cells.Item(24, 5).Value = 1.0736387843822
$ws.Cells.Item(24, 6).Value = 1.082794468117103
$ws.Cells.Item(24, 9).Value = 1.048152554080239
$ws.Cells.Item(24, 10).Value = 1.075664727866111
$ws.Cells.Item(24, 11).Value = 1.074130284487553
$ws.Cells.Item(24, 12).Value = 1.076882548718497
$ws.Cells.Item(24, 13).Value = 1.086008964048373
$ws.Cells.Item(24, 14).Value = 1.07719229516015
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.072467554976514
$ws.Cells.Item(25, 4).Value = 1.073052785231677
$ws.Cells.Item(25, 5).Value = 1.0760473803121
$ws.Cells.Item(25, 6).Value = 1.085194809964949
$ws.Cells.Item(25, 9).Value = 1.048835542674954
$ws.Cells.Item(25, 10).Value = 1.077848256715203
$ws.Cells.Item(25, 11).Value = 1.07599266489531
$ws.Cells.Item(25, 12).Value = 1.078978572956132
$ws.Cells.Item(25, 13).Value = 1.088099789255326
$ws.Cells.Item(25, 14).Value = 1.079378924870662
